{"js": "// Update the division problems in the practice-sheet table.\n// Each cell's text is a unique \"NNN\u00f7D=\" expression; replace old -> new\n// by exact text search + Replace so formatting (font/size) is preserved.\nconst replacements = [\n  [\"813\u00f72=\", \"314\u00f75=\"],\n  [\"502\u00f77=\", \"784\u00f79=\"],\n  [\"423\u00f73=\", \"957\u00f73=\"],\n  [\"445\u00f75=\", \"803\u00f76=\"],\n  [\"490\u00f76=\", \"507\u00f79=\"],\n  [\"364\u00f74=\", \"759\u00f78=\"],\n  [\"905\u00f76=\", \"571\u00f76=\"],\n  [\"305\u00f79=\", \"557\u00f77=\"],\n  [\"184\u00f76=\", \"495\u00f75=\"],\n  [\"812\u00f77=\", \"158\u00f73=\"],\n  [\"472\u00f74=\", \"933\u00f77=\"],\n  [\"403\u00f75=\", \"879\u00f78=\"],\n  [\"457\u00f78=\", \"527\u00f77=\"],\n  [\"378\u00f79=\", \"224\u00f79=\"],\n  [\"966\u00f74=\", \"269\u00f79=\"],\n  [\"970\u00f77=\", \"821\u00f79=\"],\n  [\"479\u00f77=\", \"174\u00f72=\"],\n  [\"823\u00f74=\", \"481\u00f73=\"],\n  [\"176\u00f73=\", \"833\u00f72=\"],\n  [\"452\u00f77=\", \"589\u00f74=\"],\n  [\"126\u00f77=\", \"973\u00f76=\"],\n  [\"160\u00f75=\", \"736\u00f73=\"],\n  [\"335\u00f72=\", \"369\u00f77=\"],\n  [\"222\u00f72=\", \"724\u00f73=\"],\n  [\"221\u00f75=\", \"780\u00f78=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the practice-sheet table.\n# Each cell's text is a unique \"NNN\u00f7D=\" expression; replace old -> new\n# using Find/Replace on the document body so formatting (font/size) is\n# preserved (only the text run content changes).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"813\u00f72=\"; New = \"314\u00f75=\" }\n    @{ Old = \"502\u00f77=\"; New = \"784\u00f79=\" }\n    @{ Old = \"423\u00f73=\"; New = \"957\u00f73=\" }\n    @{ Old = \"445\u00f75=\"; New = \"803\u00f76=\" }\n    @{ Old = \"490\u00f76=\"; New = \"507\u00f79=\" }\n    @{ Old = \"364\u00f74=\"; New = \"759\u00f78=\" }\n    @{ Old = \"905\u00f76=\"; New = \"571\u00f76=\" }\n    @{ Old = \"305\u00f79=\"; New = \"557\u00f77=\" }\n    @{ Old = \"184\u00f76=\"; New = \"495\u00f75=\" }\n    @{ Old = \"812\u00f77=\"; New = \"158\u00f73=\" }\n    @{ Old = \"472\u00f74=\"; New = \"933\u00f77=\" }\n    @{ Old = \"403\u00f75=\"; New = \"879\u00f78=\" }\n    @{ Old = \"457\u00f78=\"; New = \"527\u00f77=\" }\n    @{ Old = \"378\u00f79=\"; New = \"224\u00f79=\" }\n    @{ Old = \"966\u00f74=\"; New = \"269\u00f79=\" }\n    @{ Old = \"970\u00f77=\"; New = \"821\u00f79=\" }\n    @{ Old = \"479\u00f77=\"; New = \"174\u00f72=\" }\n    @{ Old = \"823\u00f74=\"; New = \"481\u00f73=\" }\n    @{ Old = \"176\u00f73=\"; New = \"833\u00f72=\" }\n    @{ Old = \"452\u00f77=\"; New = \"589\u00f74=\" }\n    @{ Old = \"126\u00f77=\"; New = \"973\u00f76=\" }\n    @{ Old = \"160\u00f75=\"; New = \"736\u00f73=\" }\n    @{ Old = \"335\u00f72=\"; New = \"369\u00f77=\" }\n    @{ Old = \"222\u00f72=\"; New = \"724\u00f73=\" }\n    @{ Old = \"221\u00f75=\"; New = \"780\u00f78=\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
